$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff. D-column cells that look like plain numbers are
# forced to remain text (matching the original inlineStr type) by briefly
# switching to a text number format and then restoring the original style,
# which avoids Excel's automatic string-to-number coercion.

$ws.Range("D2").Value = '72.298.06'
$ws.Range("E2").Value = '  +3.30%  '
$ws.Range("D3").Value = '4.028.63'
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  +0.18%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '540.29'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +2.80%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.94'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +5.75%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.699'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  +13.92%  '
$ws.Range("E8").Value = '  +0.05%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.753'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +3.50%  '
$ws.Range("E10").Value = '  +0.49%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000326'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -4.86%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.86'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +13.31%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.79'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").Value = '4.671.15'
$ws.Range("E14").Value = '  +3.30%  '
$ws.Range("D15").Value = '4.014.38'
$ws.Range("E15").Value = '  +2.72%  '
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.18'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  +0.06%  '
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.56'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  -3.61%  '
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("D20").Value = '72.083.75'
$ws.Range("E20").Value = '  +3.42%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '433.05'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +1.94%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '99.29'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +12.96%  '
$ws.Range("E23").Value = '  +2.29%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.32'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +6.96%  '
$ws.Range("E25").Value = '  +3.13%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.15'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -6.06%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.89'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +35.57%  '
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.92'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  +2.97%  '
$ws.Range("E29").Value = '  +2.48%  '
$ws.Range("E30").Value = '  +2.38%  '
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("E32").Value = '  +3.79%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '678.43'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -1.21%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.86'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -2.70%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '66.09'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -1.70%  '
$ws.Range("E36").Value = '  +6.53%  '
$ws.Range("E37").Value = '  -2.92%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.155'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +4.45%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.52'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +11.13%  '
$ws.Range("D40").Value = '0.0₃0830'
$ws.Range("E40").Value = '  -5.87%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.46'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  +1.45%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -0.06%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -0.07%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0490'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +2.06%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.152'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +6.08%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.62'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -7.22%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.60'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +7.73%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.36'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -5.13%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.03'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -3.01%  '
$ws.Range("E50").Value = '  -0.06%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '144.66'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +0.38%  '
